# Rows 5, 6 and 7 of the sheet each describe one observation record. The
# records were re-ordered (a 3-way cyclic rotation): the data that used to
# live in row 5 now lives in row 6, row 6's data moved to row 7, and row 7's
# data moved up to row 5. Columns shared by all three rows (location info,
# dates, observer names, etc.) are identical across the rows and therefore
# need no change - only the columns that actually differ are rewritten here.
#
# Note: this interop's `Range.Value` *getter* is unreliable, so reads use
# `.Value2` (numeric columns, to keep full double precision) or `.Text`
# (string columns). `Range.Value` as a *setter* works normally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the "old" values for the columns that differ between rows ---
$oldA5  = $ws.Range("A5").Value2
$oldB5  = $ws.Range("B5").Value2
$oldD5  = $ws.Range("D5").Text
$oldE5  = $ws.Range("E5").Value2
$oldF5  = $ws.Range("F5").Text
$oldG5  = $ws.Range("G5").Text
$oldH5  = $ws.Range("H5").Text
$oldI5  = $ws.Range("I5").Text
$oldQ5  = $ws.Range("Q5").Value2
$oldR5  = $ws.Range("R5").Value2
$oldAC5 = $ws.Range("AC5").Text

$oldA6  = $ws.Range("A6").Value2
$oldB6  = $ws.Range("B6").Value2
$oldD6  = $ws.Range("D6").Text
$oldE6  = $ws.Range("E6").Value2
$oldF6  = $ws.Range("F6").Text
$oldG6  = $ws.Range("G6").Text
$oldH6  = $ws.Range("H6").Text
$oldI6  = $ws.Range("I6").Text
$oldQ6  = $ws.Range("Q6").Value2
$oldR6  = $ws.Range("R6").Value2
$oldAC6 = $ws.Range("AC6").Text

$oldA7  = $ws.Range("A7").Value2
$oldB7  = $ws.Range("B7").Value2
$oldD7  = $ws.Range("D7").Text
$oldE7  = $ws.Range("E7").Value2
$oldF7  = $ws.Range("F7").Text
$oldG7  = $ws.Range("G7").Text
$oldH7  = $ws.Range("H7").Text
$oldI7  = $ws.Range("I7").Text
$oldQ7  = $ws.Range("Q7").Value2
$oldR7  = $ws.Range("R7").Value2
$oldAC7 = $ws.Range("AC7").Text

# --- row 5 <= old row 7 ---
$ws.Range("A5").Value = $oldA7
$ws.Range("B5").Value = $oldB7
$ws.Range("D5").Value = $oldD7
$ws.Range("E5").Value = $oldE7
$ws.Range("F5").Value = $oldF7
$ws.Range("G5").Value = $oldG7
$ws.Range("H5").Value = $oldH7
$ws.Range("I5").Value = "'" + $oldI7
$ws.Range("Q5").Value = $oldQ7
$ws.Range("R5").Value = $oldR7
if ($oldAC7 -eq "") {
    $ws.Range("AC5").ClearContents()
} else {
    $ws.Range("AC5").Value = $oldAC7
}

# --- row 6 <= old row 5 ---
$ws.Range("A6").Value = $oldA5
$ws.Range("B6").Value = $oldB5
$ws.Range("D6").Value = $oldD5
$ws.Range("E6").Value = $oldE5
$ws.Range("F6").Value = $oldF5
$ws.Range("G6").Value = $oldG5
$ws.Range("H6").Value = $oldH5
$ws.Range("I6").Value = "'" + $oldI5
$ws.Range("Q6").Value = $oldQ5
$ws.Range("R6").Value = $oldR5
if ($oldAC5 -eq "") {
    $ws.Range("AC6").ClearContents()
} else {
    $ws.Range("AC6").Value = $oldAC5
}

# --- row 7 <= old row 6 ---
$ws.Range("A7").Value = $oldA6
$ws.Range("B7").Value = $oldB6
$ws.Range("D7").Value = $oldD6
$ws.Range("E7").Value = $oldE6
$ws.Range("F7").Value = $oldF6
$ws.Range("G7").Value = $oldG6
$ws.Range("H7").Value = $oldH6
$ws.Range("I7").Value = "'" + $oldI6
$ws.Range("Q7").Value = $oldQ6
$ws.Range("R7").Value = $oldR6
if ($oldAC6 -eq "") {
    $ws.Range("AC7").ClearContents()
} else {
    $ws.Range("AC7").Value = $oldAC6
}
